$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.43 = 21727.97 pesos`n✅ 21727.97 pesos = 5.43 = 960.38 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the numeric rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 184.152
$ws2.Range("O10").Value = 4001.25
$ws2.Range("N12").Value = 4000
$ws2.Range("O12").Value = 176.8
